# Revert "Powerpoint writer: consolidate text run nodes."
#
# The writer used to merge a "word " run with the following run; this
# splits those merged runs back into a "word" run and a separate " "
# (space) run, matching the pre-consolidation output, without touching
# any other formatting.

$p = $ppt.ActivePresentation

# --- Slide 1: Title "Header with inline code" -----------------------
# "Header " + "with " + "inline code" (Courier)
#   -> "Header" + " " + "with" + " " + "inline code" (Courier)
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
# Split rightmost run first so earlier character offsets stay valid.
$tr1.Characters(8, 4).Text = "with"     # "with " -> "with" + " "
$tr1.Characters(1, 6).Text = "Header"   # "Header " -> "Header" + " "

# --- Slide 2: Title "Syntax highlighting" ----------------------------
# "Syntax " + "highlighting" -> "Syntax" + " " + "highlighting"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, 6).Text = "Syntax"   # "Syntax " -> "Syntax" + " "

# --- Slide 3: Title "Two column slide" -------------------------------
# "Two " + "column " + "slide" -> "Two" + " " + "column" + " " + "slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(5, 6).Text = "column"   # "column " -> "column" + " "
$tr3.Characters(1, 3).Text = "Two"      # "Two " -> "Two" + " "
